# Auto-generated PowerShell-style Excel COM-interop edit script
# Updates numeric market/profit data cells on multiple sheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 519.25
$ws.Range("I101").Value = 493.83334
$ws.Range("J101").Value = 595.5
$ws.Range("K101").Value = 1481.50002
$ws.Range("L101").Value = 1786.5
$ws.Range("M101").Value = 140.4999800000001
$ws.Range("N101").Value = -5030.5
$ws.Range("H111").Value = 371.4
$ws.Range("I111").Value = 371.4
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1114.2
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 1952.8
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1054.7778
$ws.Range("I2").Value = 999.125
$ws.Range("K2").Value = 999.125
$ws.Range("M2").Value = -886.125
$ws.Range("H92").Value = 29249.5
$ws.Range("J92").Value = 29249.5
$ws.Range("L92").Value = 29249.5
$ws.Range("N92").Value = -34241.5
$ws.Range("H116").Value = 1054.7778
$ws.Range("I116").Value = 999.125
$ws.Range("K116").Value = 999.125
$ws.Range("M116").Value = 1294.875
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1054.7778
$ws.Range("I3").Value = 999.125
$ws.Range("K3").Value = 999.125
$ws.Range("M3").Value = -885.125
$ws.Range("H92").Value = 23000
$ws.Range("J92").Value = 23000
$ws.Range("L92").Value = 23000
$ws.Range("N92").Value = -27992
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1951
$ws.Range("I6").Value = 300
$ws.Range("K6").Value = 300
$ws.Range("M6").Value = -187
$ws.Range("H58").Value = 3230
$ws.Range("I58").Value = 2699
$ws.Range("K58").Value = 2699
$ws.Range("M58").Value = -2496
$ws.Range("H86").Value = 7407.8335
$ws.Range("I86").Value = 7773.3335
$ws.Range("J86").Value = 7042.3335
$ws.Range("K86").Value = 7773.3335
$ws.Range("L86").Value = 7042.3335
$ws.Range("M86").Value = -6650.3335
$ws.Range("N86").Value = -9288.333500000001
$ws.Range("H89").Value = 7407.8335
$ws.Range("I89").Value = 7773.3335
$ws.Range("J89").Value = 7042.3335
$ws.Range("K89").Value = 38866.6675
$ws.Range("L89").Value = 35211.6675
$ws.Range("M89").Value = -33250.6675
$ws.Range("N89").Value = -46443.6675
$ws.Range("H99").Value = 5500
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -2502
$ws.Range("N99").Value = -8996
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("H126").Value = 5500
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -22940
$ws.Range("H136").Value = 3230
$ws.Range("I136").Value = 2699
$ws.Range("K136").Value = 8097
$ws.Range("M136").Value = -5547

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 146.6
$ws.Range("J12").Value = 175.75
$ws.Range("L12").Value = 527.25
$ws.Range("N12").Value = -873.25
$ws.Range("H26").Value = 85037.75
$ws.Range("I26").Value = 167549.83
$ws.Range("J26").Value = 2525.6667
$ws.Range("K26").Value = 502649.49
$ws.Range("L26").Value = 7577.000100000001
$ws.Range("M26").Value = -502361.49
$ws.Range("N26").Value = -8153.000100000001
$ws.Range("H37").Value = 99672.625
$ws.Range("J37").Value = 99672.625
$ws.Range("L37").Value = 299017.875
$ws.Range("N37").Value = -299241.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2249.1667
$ws.Range("I80").Value = 2247.5
$ws.Range("J80").Value = 2252.5
$ws.Range("K80").Value = 2247.5
$ws.Range("L80").Value = 2252.5
$ws.Range("M80").Value = -1249.5
$ws.Range("N80").Value = -4248.5
$ws.Range("H83").Value = 2249.1667
$ws.Range("I83").Value = 2247.5
$ws.Range("J83").Value = 2252.5
$ws.Range("K83").Value = 11237.5
$ws.Range("L83").Value = 11262.5
$ws.Range("M83").Value = -6245.5
$ws.Range("N83").Value = -21246.5
$ws.Range("H102").Value = 2367.9333
$ws.Range("I102").Value = 1424.7693
$ws.Range("K102").Value = 1424.7693
$ws.Range("M102").Value = 197.2307000000001
$ws.Range("H126").Value = 7398.4
$ws.Range("I126").Value = 7398.4
$ws.Range("K126").Value = 22195.2
$ws.Range("M126").Value = -19725.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3500
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H61").Value = 6995.8
$ws.Range("J61").Value = 8328
$ws.Range("L61").Value = 8328
$ws.Range("N61").Value = -8732
$ws.Range("H113").Value = 6995.8
$ws.Range("J113").Value = 8328
$ws.Range("L113").Value = 8328
$ws.Range("N113").Value = -12668
$ws.Range("H140").Value = 29476.334
$ws.Range("I140").Value = 10000
$ws.Range("J140").Value = 39214.5
$ws.Range("K140").Value = 10000
$ws.Range("L140").Value = 39214.5
$ws.Range("M140").Value = -4820
$ws.Range("N140").Value = -49574.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 133.33333
$ws.Range("I107").Value = 125
$ws.Range("K107").Value = 375
$ws.Range("M107").Value = 1545
$ws.Range("H110").Value = 20000
$ws.Range("J110").Value = 20000
$ws.Range("L110").Value = 20000
$ws.Range("N110").Value = -28180
$ws.Range("H132").Value = 3815.5715
$ws.Range("I132").Value = 3820.8
$ws.Range("J132").Value = 3802.5
$ws.Range("K132").Value = 11462.4
$ws.Range("L132").Value = 11407.5
$ws.Range("M132").Value = -8932.400000000001
$ws.Range("N132").Value = -16467.5
$ws.Range("H136").Value = 3255.6206
$ws.Range("I136").Value = 3212.2307
$ws.Range("K136").Value = 9636.6921
$ws.Range("M136").Value = -7086.6921

